# Re-ran the correlation table after excluding participants who reported
# having more than 4 majors. A few descriptive stats in the table shift
# slightly as a result:
#   - Row 3 ("2. interest.all"): M   6.03 -> 6.02
#   - Row 5 ("4. dd.id.all"):    M   6.26 -> 6.27
#                                 SD  0.82 -> 0.81
#                                 X1  .65** -> .66**
#
# The source workbook stores these figures as *text* (shared strings),
# even though they look numeric, so a straight `.Value = "6.02"` would
# make Excel auto-detect a number and reformat the cell. To keep the
# cells typed/styled exactly as before, stage each new value in a
# scratch cell that is explicitly formatted as Text, then copy just the
# value into the target cell (which keeps the target's own style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Range("H1")
$helper.NumberFormat = "@"

$targets = @("B3", "B5", "C5", "D5")
$values  = @("6.02", "6.27", "0.81", ".66**")

for ($i = 0; $i -lt $targets.Length; $i++) {
    $helper.Value = $values[$i]
    $helper.Copy()
    $ws.Range($targets[$i]).PasteSpecial(-4163)  # xlPasteValues
}

$helper.Clear()
$excel.CutCopyMode = 0
